$d = $word.ActiveDocument

# --- Change 1: split "Enrollment no: 25103002042" into three runs by
# inserting "0" after "Enrollment no: 25" (turning the enrollment number
# into "250103002042"). Track Changes is toggled on for the insertion so
# the engine produces a distinct run for the newly typed text, then the
# revision is accepted so the final XML has three plain <w:r> runs (no
# <w:ins> markup), matching a normal accepted edit.
$d.TrackRevisions = $true
$insertPoint = $d.Range(41, 41)
$insertPoint.InsertBefore("0")
$d.TrackRevisions = $false
$d.Revisions.Item(1).Accept()

# --- Change 2: mark the drawing run that carries the second
# lastRenderedPageBreak (InlineShape #4, the "Picture 1" image) as
# NoProofing, i.e. add <w:noProof/> to its run properties.
$shp = $d.InlineShapes.Item(4)
$shp.Range.NoProofing = $true
